# Re-scrape update for toilet_coop_2022-09-10.xlsx:
#  - every row's timestamp (col O) moves from 07:03:42 to 20:57:55
#  - row 3's product was replaced by a newly-scraped item (id 6568452)
#  - rows 11/12, 13/14 and 26/27 swap places (scrape order changed)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ts = "2022-09-10 20:57:55"

# Row 2 - only the timestamp changes
$ws.Range("O2").Value = $ts

# Row 3 - replaced by a different product entirely
$ws.Range("A3").Value = "'6568452"
$ws.Range("B3").Value = "Super Soft Premium Mandel feucht"
$ws.Range("C3").Value = "/de/haushalt-tier/toiletten-haushaltpapier/toilettenpapier/feuchttuecher/super-soft-premium-mandel-feucht/p/6568452"
$ws.Range("D3").Value = "50ST"
$ws.Range("E3").Value = 8
$ws.Range("F3").Value = 3.5
$ws.Range("H3").Value = "'2.95"
$ws.Range("I3").Value = "0.06/1ST"
$ws.Range("J3").Value = "Preis pro 1 Stück"
$ws.Range("K3").Value = "'0.06"
$ws.Range("L3").Value = "1ST"
$ws.Range("M3").Value = "['haushalt-tier', 'toiletten-haushaltpapier', 'toilettenpapier', 'feuchttuecher']"
$ws.Range("N3").Value = "Super Soft Premium Mandel feucht 2.95 Schweizer Franken"
$ws.Range("O3").Value = $ts

# Rows 4-10 - only the timestamp changes
foreach ($r in 4..10) {
    $ws.Range("O$r").Value = $ts
}

# Rows 11 & 12 swap places
$ws.Range("A11").Value = "'3180824"
$ws.Range("B11").Value = "Tempo Taschentücher Plus Aloe &amp; Kamille 12x9 Stück"
$ws.Range("C11").Value = "/de/haushalt-tier/toiletten-haushaltpapier/papiertaschentuecher/taschentuecher/tempo-taschentuecher-plus-aloe-kamille-12x9-stueck/p/3180824"
$ws.Range("D11").Value = "12ST"
$ws.Range("E11").Value = 7
$ws.Range("F11").Value = 4
$ws.Range("I11").Value = "0.33/1ST"
$ws.Range("K11").Value = "'0.33"
$ws.Range("N11").Value = "Tempo Taschentücher Plus Aloe &amp; Kamille 12x9 Stück 3.95 Schweizer Franken"
$ws.Range("O11").Value = $ts

$ws.Range("A12").Value = "'3754631"
$ws.Range("B12").Value = "Tempo Taschentücher sanft &amp; frei 10x9 Stück"
$ws.Range("C12").Value = "/de/haushalt-tier/toiletten-haushaltpapier/papiertaschentuecher/taschentuecher/tempo-taschentuecher-sanft-frei-10x9-stueck/p/3754631"
$ws.Range("D12").Value = "10ST"
$ws.Range("E12").Value = 9
$ws.Range("F12").Value = 3.5
$ws.Range("I12").Value = "0.40/1ST"
$ws.Range("K12").Value = "'0.40"
$ws.Range("N12").Value = "Tempo Taschentücher sanft &amp; frei 10x9 Stück 3.95 Schweizer Franken"
$ws.Range("O12").Value = $ts

# Rows 13 & 14 swap places
$ws.Range("A13").Value = "'6638996"
$ws.Range("B13").Value = "Oecoplan Allzweck Papiertücher"
$ws.Range("C13").Value = "/de/haushalt-tier/toiletten-haushaltpapier/haushaltspapier/oecoplan-allzweck-papiertuecher/p/6638996"
$ws.Range("D13").Value = "176BLT"
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 5
$ws.Range("G13").Value = "Coop"
$ws.Range("H13").Value = "'3.20"
$ws.Range("N13").Value = "Oecoplan Allzweck Papiertücher 3.20 Schweizer Franken"
$ws.Range("O13").Value = $ts

$ws.Range("A14").Value = "'6433417"
$ws.Range("B14").Value = "Plenty Fun Design extra Long"
$ws.Range("C14").Value = "/de/haushalt-tier/toiletten-haushaltpapier/haushaltspapier/plenty-fun-design-extra-long/p/6433417"
$ws.Range("D14").Value = "144BLT"
$ws.Range("E14").Value = 6
$ws.Range("F14").Value = 4.5
$ws.Range("G14").Value = "Plenty"
$ws.Range("H14").Value = "'6.70"
$ws.Range("N14").Value = "Plenty Fun Design extra Long 6.70 Schweizer Franken"
$ws.Range("O14").Value = $ts

# Rows 15-25 - only the timestamp changes
foreach ($r in 15..25) {
    $ws.Range("O$r").Value = $ts
}

# Rows 26 & 27 swap places
$ws.Range("A26").Value = "'6996129"
$ws.Range("B26").Value = "Tela Toilettenpapier Futura 3-lagig 9 Rollen"
$ws.Range("C26").Value = "/de/haushalt-tier/toiletten-haushaltpapier/toilettenpapier/toilettenpapier/tela-toilettenpapier-futura-3-lagig-9-rollen/p/6996129"
$ws.Range("D26").Value = "9Rol"
$ws.Range("E26").Value = ""
$ws.Range("F26").Value = 0
$ws.Range("G26").Value = "Tela"
$ws.Range("H26").Value = "'8.40"
$ws.Range("I26").Value = "0.93/1Rol"
$ws.Range("J26").Value = "Preis pro 1 Rolle"
$ws.Range("K26").Value = "'0.93"
$ws.Range("L26").Value = "1Rol"
$ws.Range("M26").Value = "['haushalt-tier', 'toiletten-haushaltpapier', 'toilettenpapier', 'toilettenpapier']"
$ws.Range("N26").Value = "Tela Toilettenpapier Futura 3-lagig 9 Rollen 8.40 Schweizer Franken"
$ws.Range("O26").Value = $ts

$ws.Range("A27").Value = "'6868354"
$ws.Range("B27").Value = "Tempo Bamboo Eco"
$ws.Range("C27").Value = "/de/haushalt-tier/toiletten-haushaltpapier/papiertaschentuecher/taschentuecher/tempo-bamboo-eco/p/6868354"
$ws.Range("D27").Value = "12ST"
$ws.Range("E27").Value = 1
$ws.Range("F27").Value = 3
$ws.Range("G27").Value = "Tempo"
$ws.Range("H27").Value = "'3.95"
$ws.Range("I27").Value = "0.33/1ST"
$ws.Range("J27").Value = "Preis pro 1 Stück"
$ws.Range("K27").Value = "'0.33"
$ws.Range("L27").Value = "1ST"
$ws.Range("M27").Value = "['haushalt-tier', 'toiletten-haushaltpapier', 'papiertaschentuecher', 'taschentuecher']"
$ws.Range("N27").Value = "Tempo Bamboo Eco 3.95 Schweizer Franken"
$ws.Range("O27").Value = $ts

# Rows 28-30 - only the timestamp changes
foreach ($r in 28..30) {
    $ws.Range("O$r").Value = $ts
}
